# Weekly price update: insert two new rows (new week) right after row 423,
# pushing the existing rows 424..516 down to 426..518, and populate the
# two freshly inserted rows (424, 425) with the new week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 424 (shifts 424:516 -> 426:518).
$ws.Rows("424:425").Insert()

# New "Primera" quality row for the new week.
$ws.Range("A424").Value = 8
$ws.Range("B424").Value = "Terminal La Palmera de La Serena"
$ws.Range("C424").Value = "Coquimbo"
$ws.Range("D424").Value = 45204
$ws.Range("E424").Value = 4
$ws.Range("F424").Value = 100114014
$ws.Range("G424").Value = "Betarraga"
$ws.Range("H424").Value = "Sin especificar"
$ws.Range("I424").Value = "Primera"
$ws.Range("J424").Value = 1800
$ws.Range("K424").Value = 500
$ws.Range("L424").Value = 600
$ws.Range("M424").Value = 550
$ws.Range("N424").Value = "`$/paquete 3 unidades"
$ws.Range("O424").Value = "Provincia del Elqu$([char]0xED)"
$ws.Range("P424").Value = 183
$ws.Range("Q424").Value = 3
$ws.Range("R424").Value = "Hortaliza"

# New "Segunda" quality row for the new week.
$ws.Range("A425").Value = 8
$ws.Range("B425").Value = "Terminal La Palmera de La Serena"
$ws.Range("C425").Value = "Coquimbo"
$ws.Range("D425").Value = 45204
$ws.Range("E425").Value = 4
$ws.Range("F425").Value = 100114014
$ws.Range("G425").Value = "Betarraga"
$ws.Range("H425").Value = "Sin especificar"
$ws.Range("I425").Value = "Segunda"
$ws.Range("J425").Value = 1200
$ws.Range("K425").Value = 400
$ws.Range("L425").Value = 450
$ws.Range("M425").Value = 425
$ws.Range("N425").Value = "`$/paquete 3 unidades"
$ws.Range("O425").Value = "Provincia del Elqu$([char]0xED)"
$ws.Range("P425").Value = 142
$ws.Range("Q425").Value = 3
$ws.Range("R425").Value = "Hortaliza"
